# Apply data updates to the ZoomInfo Lead Quality Analysis Report worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / summary updates -------------------------------------------------
$ws.Range("A2").Value = "Generated: 2025-07-03 09:18:31"
$ws.Range("B9").Value = 65
$ws.Range("B17").Value = "35.27s"

# --- Row 21 (Daniel J Edelman Holdings) ---------------------------------------
$ws.Range("N21").Value = @"
✅ ZI_Company_Name__c matches email domain; high brand consistency.
✅ ZI_Employees__c aligns with LS_Company_Size_Range__c; no major discrepancy.
⚠️ ZI_Website__c missing; inferred primary site from email domain.
⚠️ Large-company completeness check: website missing despite high employee count.
"@

# --- Row 22 (Medtronic) --------------------------------------------------------
$ws.Range("N22").Value = @"
⚠️ ZI_Company_Name__c 'Medtronic' does not match email domain 'emilydavisconsulting.com'.
✅ ZI_Employees__c aligns with 'Enterprise' segment size expectations.
⚠️ Website missing in ZoomInfo enrichment; inferred primary site from email domain.
"@

# --- Row 23 (Hatch) -------------------------------------------------------------
$ws.Range("N23").Value = @"
✅ ZI_Employees__c consistent with Enterprise segment size range.
⚠️ Website field blank; inferred primary site from ZoomInfo enrichment.
✅ ZI_Company_Name__c aligns with email domain, strengthening company match.
"@
$ws.Range("O23").Value = ""

# --- Row 24 (H&R Block / Heartland Business Services) ---------------------------
$ws.Range("N24").Value = @"
✅ ZI_Employees__c aligns with Enterprise segment size; no major discrepancy.
⚠️ Website field blank; inferred primary site from ZoomInfo enrichment.
"@

# --- Row 26 (Solar Turbines) -----------------------------------------------------
$ws.Range("M21").Copy()
$ws.Range("M26").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("M26").Value = 75
$ws.Range("N26").Value = @"
✅ ZI_Company_Name__c matches email domain and SegmentName 'Enterprise'.
✅ ZI_Employees__c aligns with LS_Company_Size_Range__c '1000+'.
⚠️ Website inconsistency: ZI_Website__c missing, should match Website 'www.caterpillar.com'.
"@

# --- Row 27 (Kaplan) --------------------------------------------------------------
$ws.Range("M21").Copy()
$ws.Range("M27").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("M27").Value = 75
$ws.Range("N27").Value = @"
✅ ZI_Company_Name__c matches email domain and Website.
✅ ZI_Employees__c aligns with LS_Company_Size_Range__c for Enterprise segment.
⚠️ ZI_Website__c missing; inferred as 'kaplan.com' based on email domain.
⚠️ Large-company completeness check: ZI_Website__c is not populated.
"@

# --- Row 28 (Legacy Advisor / Ameriprise Financial) --------------------------------
$ws.Range("N28").Value = @"
⚠️ ZI_Website__c (ameriprise.com) and ZI_Company_Name__c (Ameriprise Financial) do not match lead's email domain (legacyadvisor.net).
✅ ZI_Employees__c (12374) aligns with the Enterprise segment, but LS_Company_Size_Range__c is missing for comparison.
"@
$ws.Range("O28").Value = @"
{
  "ZI_Company_Name__c": "Legacy Advisor"
}
"@
$ws.Range("P28").Value = @"
{
  "ZI_Employees__c": 5000
}
"@

# --- Row 29 (Microsoft / Shields Legal Group) ----------------------------------------
$ws.Range("N29").Value = @"
❌ Company name 'Microsoft' conflicts with email domain 'shieldslegal.com'.
⚠️ Employee count (210,842) significantly exceeds Enterprise segment size expectations.
⚠️ Website inconsistency between lead-provided and enriched data.
"@
$ws.Range("P29").Value = @"
{
  "ZI_Employees__c": "100-250"
}
"@

# --- Row 30 (H&R Block via yahoo.com) --------------------------------------------------
$ws.Range("M30").Value = 50
$ws.Range("N30").Value = @"
❌ Large discrepancy in employee count (46700) and LS_Company_Size_Range__c (1000+).
⚠️ Lead has a free email domain (yahoo.com) but enriched company is H&R Block, raising authenticity concerns.
"@
